# Generate Report for Handback
# Updates the localization-status workbook to reflect that the de-de / zh-cn
# handback packages have now been generated and the overview status is
# "in sync" again.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$mdFileName  = "ba8290e7-589c-44af-9ee0-da02b12931cf.md"
$mdUrl       = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0382a5e66ccf904c9adfee5a70403dae194ae7d3/e2e/ba8290e7-589c-44af-9ee0-da02b12931cf.md"
$zhcnXlf     = "ba8290e7-589c-44af-9ee0-da02b12931cf.e0091d0906d7d380ebfa7bc3b0e53fb7c5a97800.zh-cn.xlf"
$dedeXlf     = "ba8290e7-589c-44af-9ee0-da02b12931cf.e0091d0906d7d380ebfa7bc3b0e53fb7c5a97800.de-de.xlf"

$newStatus       = "Handed back: in sync with en-US"
$zhcnHandbackDt  = "2016-08-21 05:04:29"
$dedeHandbackDt  = "2016-08-21 05:04:35"

# --- Overview sheet: status text + wider status columns -------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# --- zh-cn sheet ------------------------------------------------------------
$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

$zhcn.Range("C2").Value = $newStatus

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl, "", "", $mdFileName)
$zhcn.Range("I2").Font.Underline = 2
$zhcn.Range("I2").Font.Color = 15570276
$zhcn.Range("J2").Value = $zhcnXlf
$zhcn.Range("K2").Value = $zhcnHandbackDt

# --- de-de sheet ------------------------------------------------------------
$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40

$dede.Range("C2").Value = $newStatus

$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl, "", "", $mdFileName)
$dede.Range("I2").Font.Underline = 2
$dede.Range("I2").Font.Color = 15570276
$dede.Range("J2").Value = $dedeXlf
$dede.Range("K2").Value = $dedeHandbackDt
